# Update cryptocurrency price/volume data (and a few reordered rows)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.039.16'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.65%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.841.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.87%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.48%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '278.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.49%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.42%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5106'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.73%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3499'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.85%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.86'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.55%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06817'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8070'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.65%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07775'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.837.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.57%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.082'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.16%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9987'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.65%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.99%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008055'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9996'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.30%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.079.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.67%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.779'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.46%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.42%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.212'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.88%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.370'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.21%  '

$ws.Range('E27').Value = '  -3.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.95%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '109.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.39%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.363'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.290'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.16%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08809'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.04%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04864'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.91%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.164'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7297'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.57%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.867'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.94%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.213'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.82%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.378'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.72%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01852'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.40%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5157'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.79%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9507'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.03%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '117.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.23%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.262'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.86%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.015'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.88%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9984'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.46%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4526'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.23%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1363'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.71%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.329'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.46%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.26%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05916'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.34%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.494'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.92%  '
